# Updates symbol list data for 29-12-2022 run (Hora column 8 -> 9, refreshed prices)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.ClearFormats()
}

# Every data row in column G (Hora) moves from 8 to 9
Set-TextValue $ws.Range("G2:G51") "9"

# Row 3
Set-TextValue $ws.Range("D3") "23.90"

# Row 4
Set-TextValue $ws.Range("D4") "5.186"

# Row 5
Set-TextValue $ws.Range("D5") "0.05716"

# Row 6
Set-TextValue $ws.Range("D6") "6.484"

# Row 7
Set-TextValue $ws.Range("D7") "3.170"

# Row 8
Set-TextValue $ws.Range("D8") "0.8130"

# Row 9
Set-TextValue $ws.Range("D9") "0.8550"

# Row 10
Set-TextValue $ws.Range("D10") "0.1372"

# Row 11
Set-TextValue $ws.Range("D11") "0.06932"

# Row 12
Set-TextValue $ws.Range("D12") "0.03183"
Set-TextValue $ws.Range("E12") "11LiechtensteinCryptoassetsExchangeLCX"

# Row 14
Set-TextValue $ws.Range("D14") "0.09329"

# Row 15
Set-TextValue $ws.Range("D15") "3.815"

# Row 16
Set-TextValue $ws.Range("D16") "0.001525"

# Row 17
Set-TextValue $ws.Range("D17") "0.04702"

# Row 18
Set-TextValue $ws.Range("B18") "TigerCash"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D18") "0.006200"
Set-TextValue $ws.Range("E18") "17TigerCashTCH"

# Row 19
Set-TextValue $ws.Range("B19") "BitKan"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws.Range("D19") "0.001237"
Set-TextValue $ws.Range("E19") "18BitKanKAN"

# Row 20
Set-TextValue $ws.Range("B20") "HotbitToken"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws.Range("D20") "0.004108"
Set-TextValue $ws.Range("E20") "19HotbitTokenHTB"

# Row 21
Set-TextValue $ws.Range("B21") "NitroEx"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws.Range("D21") "0.00008503"
Set-TextValue $ws.Range("E21") "20NitroExNTX"

# Row 22
Set-TextValue $ws.Range("B22") "LEO"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D22") "3.540"
Set-TextValue $ws.Range("E22") "21LEOLEO"

# Row 23
Set-TextValue $ws.Range("B23") "BTSEToken"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D23") "2.157"
Set-TextValue $ws.Range("E23") "22BTSETokenBTSE"

# Row 24
Set-TextValue $ws.Range("B24") "One"
Set-TextValue $ws.Range("C24") "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D24") "0.0005982"
Set-TextValue $ws.Range("E24") "23OneONE"

# Row 40
Set-TextValue $ws.Range("D40") "0.03685"

# Row 41
Set-TextValue $ws.Range("D41") "0.006383"
Set-TextValue $ws.Range("E41") "40KickTokenKICKBestin24h"

# Row 42
Set-TextValue $ws.Range("D42") "0.1051"

# Row 44
Set-TextValue $ws.Range("D44") "0.007877"

# Row 45
Set-TextValue $ws.Range("D45") "0.00005455"

# Row 48
Set-TextValue $ws.Range("D48") "0.002557"
Set-TextValue $ws.Range("E48") "47BOLOBOLOWorstin24h"
